$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 160.5
$ws.Range("I18").Value = 160.5
$ws.Range("K18").Value = 160.5
$ws.Range("M18").Value = 123.5
$ws.Range("H40").Value = 2101.6128
$ws.Range("J40").Value = 2322.7693
$ws.Range("L40").Value = 2322.7693
$ws.Range("N40").Value = -2672.7693
$ws.Range("H129").Value = 844.2222
$ws.Range("J129").Value = 901.25
$ws.Range("L129").Value = 2703.75
$ws.Range("N129").Value = -12703.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 165.85715
$ws.Range("I5").Value = 176.83333
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 176.83333
$ws.Range("L5").Value = 100
$ws.Range("M5").Value = -64.83332999999999
$ws.Range("N5").Value = -324
$ws.Range("H61").Value = 2703.6
$ws.Range("I61").Value = 2776.8333
$ws.Range("J61").Value = 2593.75
$ws.Range("K61").Value = 2776.8333
$ws.Range("L61").Value = 2593.75
$ws.Range("M61").Value = -2564.8333
$ws.Range("N61").Value = -3017.75
$ws.Range("H74").Value = 5645.852
$ws.Range("I74").Value = 7421.3335
$ws.Range("J74").Value = 2094.889
$ws.Range("K74").Value = 7421.3335
$ws.Range("L74").Value = 2094.889
$ws.Range("M74").Value = -6547.3335
$ws.Range("N74").Value = -3842.889
$ws.Range("H77").Value = 5645.852
$ws.Range("I77").Value = 7421.3335
$ws.Range("J77").Value = 2094.889
$ws.Range("K77").Value = 37106.6675
$ws.Range("L77").Value = 10474.445
$ws.Range("M77").Value = -32738.6675
$ws.Range("N77").Value = -19210.445
$ws.Range("H136").Value = 2703.6
$ws.Range("I136").Value = 2776.8333
$ws.Range("J136").Value = 2593.75
$ws.Range("K136").Value = 8330.499899999999
$ws.Range("L136").Value = 7781.25
$ws.Range("M136").Value = -5780.499899999999
$ws.Range("N136").Value = -12881.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 165.85715
$ws.Range("I4").Value = 176.83333
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 176.83333
$ws.Range("L4").Value = 100
$ws.Range("M4").Value = -61.83332999999999
$ws.Range("N4").Value = -330
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H86").Value = 3141.5454
$ws.Range("I86").Value = 3266.1875
$ws.Range("J86").Value = 2809.1667
$ws.Range("K86").Value = 3266.1875
$ws.Range("L86").Value = 2809.1667
$ws.Range("M86").Value = -2143.1875
$ws.Range("N86").Value = -5055.1667
$ws.Range("H89").Value = 3141.5454
$ws.Range("I89").Value = 3266.1875
$ws.Range("J89").Value = 2809.1667
$ws.Range("K89").Value = 16330.9375
$ws.Range("L89").Value = 14045.8335
$ws.Range("M89").Value = -10714.9375
$ws.Range("N89").Value = -25277.8335
$ws.Range("H99").Value = 1298.3043
$ws.Range("I99").Value = 975
$ws.Range("J99").Value = 1470.7333
$ws.Range("K99").Value = 975
$ws.Range("L99").Value = 1470.7333
$ws.Range("M99").Value = 523
$ws.Range("N99").Value = -4466.7333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H41").Value = 12000
$ws.Range("I41").Value = 12000
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 12000
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -11572
$ws.Range("N41").ClearContents()
$ws.Range("H48").Value = 30000
$ws.Range("J48").Value = 30000
$ws.Range("L48").Value = 30000
$ws.Range("N48").Value = -30952
$ws.Range("H50").Value = 30390.25
$ws.Range("J50").Value = 30390.25
$ws.Range("L50").Value = 30390.25
$ws.Range("N50").Value = -31640.25
$ws.Range("H51").Value = 33164.07
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 33164.07
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 33164.07
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -34636.07
$ws.Range("H59").Value = 23075
$ws.Range("J59").Value = 35812.5
$ws.Range("L59").Value = 35812.5
$ws.Range("N59").Value = -38102.5
$ws.Range("H60").Value = 26733.934
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 26733.934
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 26733.934
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -27755.934
$ws.Range("H61").Value = 33164.07
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 33164.07
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 33164.07
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -33860.07
$ws.Range("H62").Value = 9502.75
$ws.Range("J62").Value = 9335.333000000001
$ws.Range("L62").Value = 9335.333000000001
$ws.Range("N62").Value = -10583.333
$ws.Range("H65").Value = 9502.75
$ws.Range("J65").Value = 9335.333000000001
$ws.Range("L65").Value = 46676.665
$ws.Range("N65").Value = -52916.665
$ws.Range("H68").Value = 39800
$ws.Range("J68").Value = 39800
$ws.Range("L68").Value = 39800
$ws.Range("N68").Value = -41298
$ws.Range("H71").Value = 39800
$ws.Range("J71").Value = 39800
$ws.Range("L71").Value = 119400
$ws.Range("N71").Value = -126888

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 1315.3334
$ws.Range("I130").Value = 746
$ws.Range("J130").Value = 1600
$ws.Range("K130").Value = 2238
$ws.Range("L130").Value = 4800
$ws.Range("M130").Value = 2782
$ws.Range("N130").Value = -14840
$ws.Range("H131").Value = 830.3
$ws.Range("J131").Value = 830.3
$ws.Range("L131").Value = 2490.9
$ws.Range("N131").Value = -12570.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 9761.666999999999
$ws.Range("J109").Value = 9761.666999999999
$ws.Range("L109").Value = 9761.666999999999
$ws.Range("N109").Value = -11841.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 145143
$ws.Range("I46").Value = 169166.83
$ws.Range("J46").Value = 1000
$ws.Range("K46").Value = 169166.83
$ws.Range("L46").Value = 1000
$ws.Range("M46").Value = -168978.83
$ws.Range("N46").Value = -1376
$ws.Range("H68").Value = 1584.5
$ws.Range("I68").Value = 1398
$ws.Range("J68").Value = 1646.6666
$ws.Range("K68").Value = 1398
$ws.Range("L68").Value = 1646.6666
$ws.Range("M68").Value = -649
$ws.Range("N68").Value = -3144.6666
$ws.Range("H71").Value = 1584.5
$ws.Range("I71").Value = 1398
$ws.Range("J71").Value = 1646.6666
$ws.Range("K71").Value = 6990
$ws.Range("L71").Value = 8233.333000000001
$ws.Range("M71").Value = -3246
$ws.Range("N71").Value = -15721.333
$ws.Range("H136").Value = 11180547
$ws.Range("I136").Value = 136227.06
$ws.Range("J136").Value = 22224866
$ws.Range("K136").Value = 408681.18
$ws.Range("L136").Value = 66674598
$ws.Range("M136").Value = -406131.18
$ws.Range("N136").Value = -66679698
$ws.Range("H137").Value = 30429
$ws.Range("J137").Value = 30429
$ws.Range("L137").Value = 30429
$ws.Range("N137").Value = -40629

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 13510.5
$ws.Range("J33").Value = 13510.5
$ws.Range("L33").Value = 13510.5
$ws.Range("N33").Value = -14010.5
$ws.Range("H36").Value = 13510.5
$ws.Range("J36").Value = 13510.5
$ws.Range("L36").Value = 13510.5
$ws.Range("N36").Value = -14010.5

